$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("B2").Value = 13492947.6276
$ws.Range("C2").Value = 530225.1314000001
$ws.Range("D2").Value = 255848.158
$ws.Range("F2").Value = 3188

# Update existing row 3 values
$ws.Range("B3").Value = 10430611.779336
$ws.Range("C3").Value = 329870.800764
$ws.Range("D3").Value = 837904.187
$ws.Range("F3").Value = 2447

# Add new row 4
$ws.Range("A4").Value = "bs3"
$ws.Range("B4").Value = 37775685.81309601
$ws.Range("C4").Value = 1542039.631629
$ws.Range("D4").Value = 255578.218
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 3101
